$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix missing values in row 3 (D3, F3)
$ws.Range("D3").Value = 17631
$ws.Range("F3").Value = 50968

# Update row 4 (2015年) with more precise decimal values
$ws.Range("B4").Value = 29437.64
$ws.Range("C4").Value = 11894.02
$ws.Range("D4").Value = 19320.11
$ws.Range("E4").Value = 5221.17
$ws.Range("F4").Value = 54543.54

# Update row 5 (2016年)
$ws.Range("B5").Value = 31990.4
$ws.Range("C5").Value = 12898.92
$ws.Range("D5").Value = 20924.35
$ws.Range("E5").Value = 5528.7
$ws.Range("F5").Value = 59259.46

# Update row 6 (2017年)
$ws.Range("B6").Value = 34546.8408921384
$ws.Range("C6").Value = 13842.7815106234
$ws.Range("D6").Value = 22495.316792668
$ws.Range("E6").Value = 5958.426539016
$ws.Range("F6").Value = 64933.9617724008

# Update row 7 (2018年)
$ws.Range("B7").Value = 36471.36
$ws.Range("C7").Value = 14360.52
$ws.Range("D7").Value = 23188.9
$ws.Range("E7").Value = 6440.48
$ws.Range("F7").Value = 70639.50999999999

# Update row 8 (2019年)
$ws.Range("B8").Value = 39230.49
$ws.Range("C8").Value = 15776.96
$ws.Range("D8").Value = 25034.72
$ws.Range("E8").Value = 7380.38
$ws.Range("F8").Value = 76400.72

# Update row 9 (2020年)
$ws.Range("B9").Value = 41171.74
$ws.Range("C9").Value = 16442.72
$ws.Range("D9").Value = 26248.87
$ws.Range("E9").Value = 7868.75
$ws.Range("F9").Value = 80293.83

# Add new row 10 (2021年) - copy style from A9 then set values
$ws.Range("A9").Copy()
$ws.Range("A10").PasteSpecial(-4122)
$ws.Range("A10").Value = "2021年"
$ws.Range("B10").Value = 44948.9
$ws.Range("C10").Value = 18445
$ws.Range("D10").Value = 29053.3
$ws.Range("E10").Value = 8332.799999999999
$ws.Range("F10").Value = 85835.8

# Add new row 11 (2022年) - copy style from A9 then set values
$ws.Range("A9").Copy()
$ws.Range("A11").PasteSpecial(-4122)
$ws.Range("A11").Value = "2022年"
$ws.Range("B11").Value = 47397.41
$ws.Range("C11").Value = 19302.71
$ws.Range("D11").Value = 30598.3
$ws.Range("E11").Value = 8601.129999999999
$ws.Range("F11").Value = 90116.28999999999
